# Kalyankar_LabExam03Grading.xlsx — grading update for the
# "CustomerMappingDriver Class" (rows 29-30) and the compilation-errors
# row (37), per commit "from 33-41 - Driver".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 29: S.No 18, addProduct()/getProducts() driver test -----------
# Partial credit for scanning input data, with a grading comment.
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "For not properly scanning data from input file"

# --- Row 30: S.No 19, output check --------------------------------------
# No credit for missing output, with a grading comment.
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "For no output"

# --- Row 37: compilation errors / exceptions deduction ------------------
# Half of the -5 deduction applied, with a grading comment.
$ws.Range("E37").Value = -2.5
$ws.Range("F37").Value = "For getting exceptions"

# Totals in rows 31 and 38 are formula-driven (SUM) and recalculate
# automatically.

# Leave the view pointed at the last-edited cell, matching the author's
# final on-screen selection.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("F37").Select()

Write-Host "Applied driver-class (rows 29-30) and compilation-error (row 37) grading updates."
